# Generate Report for Handback
# Regenerating the handback report updates the handoff/handback timestamps
# recorded for the file that was just processed (f5348948-...), both on the
# per-locale detail sheets and on the Overview summary sheet.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 corresponds to f5348948-4f86-4d36-b8a0-67a8c1d3ffcf.md
# Column G = "Latest HO Xliff Generate Date"
$ws_overview.Range("G3").Value = "2016-10-17 13:51:16"

# zh-cn detail sheet: row 3 corresponds to f5348948-4f86-4d36-b8a0-67a8c1d3ffcf.md
# Column H = "Correspond Handoff Datetime", Column K = "Correspond Handback DateTime"
$ws_zhcn.Range("H3").Value = "2016-10-17 13:50:54"
$ws_zhcn.Range("K3").Value = "2016-10-17 13:52:00"

# de-de detail sheet: row 3 corresponds to f5348948-4f86-4d36-b8a0-67a8c1d3ffcf.md
# Column H = "Correspond Handoff Datetime", Column K = "Correspond Handback DateTime"
$ws_dede.Range("H3").Value = "2016-10-17 13:51:16"
$ws_dede.Range("K3").Value = "2016-10-17 13:52:38"
